$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Fold_1")
$ws.Range("B2").Value = 32.0124
$ws.Range("C2").Value = 6.269600000000001
$ws.Range("D2").Value = 142.6375
$ws.Range("F2").Value = 37.64915000000001
$ws.Range("G2").Value = 37.64915000000001
$ws.Range("K2").Value = 42.2058
$ws.Range("L2").Value = 37.649
$ws.Range("M2").Value = 4.5568
$ws.Range("N2").Value = 4.5568
$ws.Range("B3").Value = 61.255
$ws.Range("C3").Value = 2.449
$ws.Range("D3").Value = 197.088
$ws.Range("F3").Value = 50.383
$ws.Range("G3").Value = 47.644
$ws.Range("H3").Value = 2.739
$ws.Range("I3").Value = 2.739
$ws.Range("K3").Value = 47.644
$ws.Range("L3").Value = 47.644
$ws.Range("B4").Value = 79.155
$ws.Range("C4").Value = 3.078
$ws.Range("D4").Value = 197.085
$ws.Range("F4").Value = 53.391
$ws.Range("G4").Value = 53.39058082191782
$ws.Range("K4").Value = 53.391
$ws.Range("L4").Value = 53.391

$ws = $wb.Worksheets.Item("Fold_2")
$ws.Range("B2").Value = 26.79765
$ws.Range("C2").Value = 6.501450000000001
$ws.Range("D2").Value = 146.8535
$ws.Range("F2").Value = 36.10555000000001
$ws.Range("G2").Value = 36.10555000000001
$ws.Range("K2").Value = 76.877
$ws.Range("L2").Value = 36.106
$ws.Range("M2").Value = 40.771
$ws.Range("N2").Value = 40.771
$ws.Range("B3").Value = 37.668
$ws.Range("C3").Value = 10.05
$ws.Range("D3").Value = 177.886
$ws.Range("F3").Value = 45.912
$ws.Range("G3").Value = 40.568
$ws.Range("H3").Value = 5.344
$ws.Range("I3").Value = 5.344
$ws.Range("K3").Value = 63.98599999999999
$ws.Range("L3").Value = 40.568
$ws.Range("M3").Value = 23.418
$ws.Range("N3").Value = 23.418
$ws.Range("B4").Value = 70.002
$ws.Range("C4").Value = 7.605
$ws.Range("D4").Value = 187.794
$ws.Range("F4").Value = 50.774
$ws.Range("G4").Value = 50.7739095890411
$ws.Range("K4").Value = 53.6682
$ws.Range("L4").Value = 50.774
$ws.Range("M4").Value = 2.8942
$ws.Range("N4").Value = 2.8942

$ws = $wb.Worksheets.Item("Fold_3")
$ws.Range("B2").Value = 31.1062
$ws.Range("C2").Value = 4.577599999999999
$ws.Range("D2").Value = 145.9664
$ws.Range("F2").Value = 37.25110000000001
$ws.Range("G2").Value = 37.25110000000001
$ws.Range("K2").Value = 50.1186
$ws.Range("L2").Value = 37.251
$ws.Range("M2").Value = 12.8676
$ws.Range("N2").Value = 12.8676
$ws.Range("B3").Value = 61.255
$ws.Range("C3").Value = 2.449
$ws.Range("D3").Value = 197.088
$ws.Range("F3").Value = 50.383
$ws.Range("G3").Value = 47.644
$ws.Range("H3").Value = 2.739
$ws.Range("I3").Value = 2.739
$ws.Range("K3").Value = 47.644
$ws.Range("L3").Value = 47.644
$ws.Range("B4").Value = 79.155
$ws.Range("C4").Value = 3.078
$ws.Range("D4").Value = 197.085
$ws.Range("F4").Value = 53.391
$ws.Range("G4").Value = 53.39058082191782
$ws.Range("K4").Value = 53.391
$ws.Range("L4").Value = 53.391

$ws = $wb.Worksheets.Item("Fold_4")
$ws.Range("B2").Value = 29.25955
$ws.Range("C2").Value = 5.13815
$ws.Range("D2").Value = 140.6869
$ws.Range("F2").Value = 36.6158
$ws.Range("G2").Value = 36.6158
$ws.Range("K2").Value = 66.73939999999999
$ws.Range("L2").Value = 36.616
$ws.Range("M2").Value = 30.1234
$ws.Range("N2").Value = 30.1234
$ws.Range("B3").Value = 40.769
$ws.Range("C3").Value = 14.381
$ws.Range("D3").Value = 172.599
$ws.Range("F3").Value = 47.145
$ws.Range("G3").Value = 41.902
$ws.Range("H3").Value = 5.243
$ws.Range("I3").Value = 5.243
$ws.Range("K3").Value = 58.37579999999999
$ws.Range("L3").Value = 41.902
$ws.Range("M3").Value = 16.4738
$ws.Range("N3").Value = 16.4738
$ws.Range("B4").Value = 79.155
$ws.Range("C4").Value = 3.078
$ws.Range("D4").Value = 195.798
$ws.Range("F4").Value = 53.364
$ws.Range("G4").Value = 53.36411780821918
$ws.Range("K4").Value = 53.4744
$ws.Range("L4").Value = 53.364
$ws.Range("M4").Value = 0.1104
$ws.Range("N4").Value = 0.1104

$ws = $wb.Worksheets.Item("Fold_5")
$ws.Range("B2").Value = 31.9586
$ws.Range("C2").Value = 5.421200000000001
$ws.Range("D2").Value = 137.5753
$ws.Range("F2").Value = 37.4382
$ws.Range("G2").Value = 37.4382
$ws.Range("K2").Value = 53.84220000000001
$ws.Range("L2").Value = 37.438
$ws.Range("M2").Value = 16.4042
$ws.Range("N2").Value = 16.4042
$ws.Range("O2").Value = 0
$ws.Range("B3").Value = 61.255
$ws.Range("C3").Value = 2.449
$ws.Range("D3").Value = 186.236
$ws.Range("F3").Value = 50.16
$ws.Range("G3").Value = 47.421
$ws.Range("H3").Value = 2.739
$ws.Range("I3").Value = 2.739
$ws.Range("K3").Value = 49.4588
$ws.Range("L3").Value = 47.421
$ws.Range("M3").Value = 2.0378
$ws.Range("O3").Value = 0
$ws.Range("B4").Value = 79.155
$ws.Range("C4").Value = 3.078
$ws.Range("D4").Value = 197.085
$ws.Range("F4").Value = 53.391
$ws.Range("G4").Value = 53.39058082191782
$ws.Range("K4").Value = 53.391
$ws.Range("L4").Value = 53.391
